# Update "想去人数" (number of people interested) counts for a handful of
# events on both the "展览" and "全部类型" worksheets, matching the data
# refresh that happened upstream (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value for each affected worksheet.
$updates = @{
    "F7"  = 1505
    "F19" = 3552
    "F21" = 317
    "F26" = 1308
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
